$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.044584707210908
$ws.Range("D2").Value = 1.051159137089707
$ws.Range("E2").Value = 1.052579708349504
$ws.Range("F2").Value = 1.063719870314751
$ws.Range("I2").Value = 1.044447491456617
$ws.Range("J2").Value = 1.049649007550051
$ws.Range("K2").Value = 1.053911412809507
$ws.Range("L2").Value = 1.055328049148805
$ws.Range("M2").Value = 1.066437747808337
$ws.Range("N2").Value = 1.051139629537196
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.045405323160177
$ws.Range("D3").Value = 1.051788138586206
$ws.Range("E3").Value = 1.053294774325047
$ws.Range("F3").Value = 1.064480036127936
$ws.Range("I3").Value = 1.044634618538296
$ws.Range("J3").Value = 1.050117576079361
$ws.Range("K3").Value = 1.054353477440171
$ws.Range("L3").Value = 1.055856240017486
$ws.Range("M3").Value = 1.067013116890383
$ws.Range("N3").Value = 1.051608863487563
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.045937104344413
$ws.Range("D4").Value = 1.05219584213004
$ws.Range("E4").Value = 1.053758529399796
$ws.Range("F4").Value = 1.064973012575787
$ws.Range("I4").Value = 1.044754943806421
$ws.Range("J4").Value = 1.050420869894402
$ws.Range("K4").Value = 1.054639527564404
$ws.Range("L4").Value = 1.05619839885712
$ws.Range("M4").Value = 1.067385849954714
$ws.Range("N4").Value = 1.051912588014609
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.046160851893282
$ws.Range("D5").Value = 1.052367405772294
$ws.Range("E5").Value = 1.053953743686296
$ws.Range("F5").Value = 1.065180520694473
$ws.Range("I5").Value = 1.044805346189354
$ws.Range("J5").Value = 1.05054839683712
$ws.Range("K5").Value = 1.054759782769839
$ws.Range("L5").Value = 1.056342332811121
$ws.Range("M5").Value = 1.067542648354779
$ws.Range("N5").Value = 1.052040296060214
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.046198430987721
$ws.Range("D6").Value = 1.052396221686969
$ws.Range("E6").Value = 1.053986535722533
$ws.Range("F6").Value = 1.065215377453508
$ws.Range("I6").Value = 1.044813798254882
$ws.Range("J6").Value = 1.050569810434616
$ws.Range("K6").Value = 1.054779974091524
$ws.Range("L6").Value = 1.056366505209782
$ws.Range("M6").Value = 1.067568981393574
$ws.Range("N6").Value = 1.052061740067475
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.04594009333868
$ws.Range("D7").Value = 1.052198133925611
$ws.Range("E7").Value = 1.053761136876351
$ws.Range("F7").Value = 1.064975784287251
$ws.Range("I7").Value = 1.044755618003051
$ws.Range("J7").Value = 1.050422573830078
$ws.Range("K7").Value = 1.054641134423154
$ws.Range("L7").Value = 1.056200321756494
$ws.Range("M7").Value = 1.067387944706252
$ws.Range("N7").Value = 1.051914294370069
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.044861873976339
$ws.Range("D8").Value = 1.051371565514956
$ws.Range("E8").Value = 1.052821147643609
$ws.Range("F8").Value = 1.063976543216177
$ws.Range("I8").Value = 1.044510888505213
$ws.Range("J8").Value = 1.049807341215963
$ws.Range("K8").Value = 1.05406080879422
$ws.Range("L8").Value = 1.055506473272832
$ws.Range("M8").Value = 1.066632106254767
$ws.Range("N8").Value = 1.051298188055074
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.042968029646782
$ws.Range("D9").Value = 1.049920473141691
$ws.Range("E9").Value = 1.051172970521424
$ws.Range("F9").Value = 1.062224258667227
$ws.Range("I9").Value = 1.044073869553238
$ws.Range("J9").Value = 1.04872403987588
$ws.Range("K9").Value = 1.053038299358119
$ws.Range("L9").Value = 1.054286833778494
$ws.Range("M9").Value = 1.065303597499149
$ws.Range("N9").Value = 1.050213348302832
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.041709690055002
$ws.Range("D10").Value = 1.048956841397775
$ws.Range("E10").Value = 1.050079822014488
$ws.Range("F10").Value = 1.061061914963809
$ws.Range("I10").Value = 1.04377868812741
$ws.Range("J10").Value = 1.04800247170155
$ws.Range("K10").Value = 1.052356776594583
$ws.Range("L10").Value = 1.05347585795932
$ws.Range("M10").Value = 1.064420296798249
$ws.Range("N10").Value = 1.049490755418929
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.041165840355169
$ws.Range("D11").Value = 1.048540495080731
$ws.Range("E11").Value = 1.049607839012735
$ws.Range("F11").Value = 1.060560020870199
$ws.Range("I11").Value = 1.043649971511835
$ws.Range("J11").Value = 1.04769019243444
$ws.Range("K11").Value = 1.052061723134002
$ws.Range("L11").Value = 1.053125218417605
$ws.Range("M11").Value = 1.064038402141622
$ws.Range("N11").Value = 1.049178032679442
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.040963985455772
$ws.Range("D12").Value = 1.048385984606943
$ws.Range("E12").Value = 1.049432729546504
$ws.Range("F12").Value = 1.060373808854667
$ws.Range("I12").Value = 1.043602025840192
$ws.Range("J12").Value = 1.047574224155289
$ws.Range("K12").Value = 1.05195213615371
$ws.Range("L12").Value = 1.052995054761019
$ws.Range("M12").Value = 1.063896638424368
$ws.Range("N12").Value = 1.049061899712033
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.04100727696798
$ws.Range("D13").Value = 1.048419121302361
$ws.Range("E13").Value = 1.049470281764447
$ws.Range("F13").Value = 1.060413742253993
$ws.Range("I13").Value = 1.043612316427238
$ws.Range("J13").Value = 1.047599098538785
$ws.Range("K13").Value = 1.051975642502029
$ws.Range("L13").Value = 1.053022971684719
$ws.Range("M13").Value = 1.063927043168865
$ws.Range("N13").Value = 1.049086809420006
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.041149151791326
$ws.Range("D14").Value = 1.048527720348012
$ws.Range("E14").Value = 1.049593360194523
$ws.Range("F14").Value = 1.060544624152777
$ws.Range("I14").Value = 1.043646011048659
$ws.Range("J14").Value = 1.047680605918701
$ws.Range("K14").Value = 1.052052664447847
$ws.Range("L14").Value = 1.053114457418462
$ws.Range("M14").Value = 1.064026682083517
$ws.Range("N14").Value = 1.049168432549751
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.041236586160905
$ws.Range("D15").Value = 1.048594650293709
$ws.Range("E15").Value = 1.049669220241417
$ws.Range("F15").Value = 1.06062529320634
$ws.Range("I15").Value = 1.043666753606287
$ws.Range("J15").Value = 1.04773082881273
$ws.Range("K15").Value = 1.052100121448018
$ws.Range("L15").Value = 1.053170835378535
$ws.Range("M15").Value = 1.06408808474851
$ws.Range("N15").Value = 1.049218726766049
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.041745805213665
$ws.Range("D16").Value = 1.048984492327895
$ws.Range("E16").Value = 1.050111174751138
$ws.Range("F16").Value = 1.061095253883192
$ws.Range("I16").Value = 1.043787211707483
$ws.Range("J16").Value = 1.048023200236099
$ws.Range("K16").Value = 1.052376359505784
$ws.Range("L16").Value = 1.053499139808895
$ws.Range("M16").Value = 1.064445654267589
$ws.Range("N16").Value = 1.049511513390375
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.042065499203401
$ws.Range("D17").Value = 1.049229275660744
$ws.Range("E17").Value = 1.050388766056383
$ws.Range("F17").Value = 1.06139042642307
$ws.Range("I17").Value = 1.043862531197699
$ws.Range("J17").Value = 1.048206642159137
$ws.Range("K17").Value = 1.052549650850547
$ws.Range("L17").Value = 1.053705216413177
$ws.Range("M17").Value = 1.064670104810959
$ws.Range("N17").Value = 1.049695215821986
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.042252069552165
$ws.Range("D18").Value = 1.049372141556152
$ws.Range("E18").Value = 1.050550811209166
$ws.Range("F18").Value = 1.061562731314339
$ws.Range("I18").Value = 1.043906376784741
$ws.Range("J18").Value = 1.048313656506668
$ws.Range("K18").Value = 1.052650733525691
$ws.Range("L18").Value = 1.053825467314469
$ws.Range("M18").Value = 1.064801078894893
$ws.Range("N18").Value = 1.049802382142165
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.042315701827287
$ws.Range("D19").Value = 1.049420869986761
$ws.Range("E19").Value = 1.050606086548082
$ws.Range("F19").Value = 1.061621505786811
$ws.Range("I19").Value = 1.043921312216527
$ws.Range("J19").Value = 1.048350148266341
$ws.Range("K19").Value = 1.052685200870181
$ws.Range("L19").Value = 1.053866478180028
$ws.Range("M19").Value = 1.064845747097176
$ws.Range("N19").Value = 1.049838925724322
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.042031188887329
$ws.Range("D20").Value = 1.049203003617606
$ws.Range("E20").Value = 1.050358969587911
$ws.Range("F20").Value = 1.061358743146652
$ws.Range("I20").Value = 1.043854459122048
$ws.Range("J20").Value = 1.04818695893464
$ws.Range("K20").Value = 1.052531057825948
$ws.Range("L20").Value = 1.05368310117038
$ws.Range("M20").Value = 1.064646017599226
$ws.Range("N20").Value = 1.049675504645052
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.041107368895332
$ws.Range("D21").Value = 1.048495736780479
$ws.Range("E21").Value = 1.049557110960196
$ws.Range("F21").Value = 1.0605060767762
$ws.Range("I21").Value = 1.043636092531426
$ws.Range("J21").Value = 1.047656603293636
$ws.Range("K21").Value = 1.052029983144172
$ws.Range("L21").Value = 1.053087514943122
$ws.Range("M21").Value = 1.063997338436612
$ws.Range("N21").Value = 1.049144395838205
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.040527424870736
$ws.Range("D22").Value = 1.048051855736535
$ws.Range("E22").Value = 1.049054143703947
$ws.Range("F22").Value = 1.059971209642323
$ws.Range("I22").Value = 1.043498018488933
$ws.Range("J22").Value = 1.047323299625704
$ws.Range("K22").Value = 1.051714990459139
$ws.Range("L22").Value = 1.052713506892939
$ws.Range("M22").Value = 1.063590003682121
$ws.Range("N22").Value = 1.048810618840834
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.040834778388964
$ws.Range("D23").Value = 1.048287088446035
$ws.Range("E23").Value = 1.049320662376599
$ws.Range("F23").Value = 1.060254634742145
$ws.Range("I23").Value = 1.043571287660183
$ws.Range("J23").Value = 1.047499975324082
$ws.Range("K23").Value = 1.05188196855017
$ws.Range("L23").Value = 1.052911731436291
$ws.Range("M23").Value = 1.063805890081589
$ws.Range("N23").Value = 1.048987545438969
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.042046691938952
$ws.Range("D24").Value = 1.049214874551647
$ws.Range("E24").Value = 1.050372432925446
$ws.Range("F24").Value = 1.06137305903635
$ws.Range("I24").Value = 1.043858106814488
$ws.Range("J24").Value = 1.048195852887975
$ws.Range("K24").Value = 1.052539459198904
$ws.Range("L24").Value = 1.053693093942679
$ws.Range("M24").Value = 1.064656901400965
$ws.Range("N24").Value = 1.049684411228821
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.043456897575487
$ws.Range("D25").Value = 1.050294960393523
$ws.Range("E25").Value = 1.051598079190699
$ws.Range("F25").Value = 1.062676245129116
$ws.Range("I25").Value = 1.044187528891325
$ws.Range("J25").Value = 1.049003993819234
$ws.Range("K25").Value = 1.053302622232869
$ws.Range("L25").Value = 1.05460177348379
$ws.Range("M25").Value = 1.065646638314294
$ws.Range("N25").Value = 1.050493699812897

Write-Output "updated cells"
